$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.003420333333333333
$ws.Range("H2").Value = 0.010261
$ws.Range("I2").Value = 0.003549653112303053
$ws.Range("J2").Value = 0.003549653112303053
$ws.Range("M2").Value = 1.424719
$ws.Range("N2").Value = 4.274157
$ws.Range("O2").Value = 0.07423298812267187
$ws.Range("P2").Value = 0.07423298812267187
$ws.Range("Q2").Value = 0.004873013886333332
$ws.Range("R2").Value = 0.04385712497699999
$ws.Range("S2").Value = 0.0002635013573251978
$ws.Range("T2").Value = 0.0002635013573251978

$ws.Range("G3").Value = 0.003420333333333333
$ws.Range("H3").Value = 0.010261
$ws.Range("I3").Value = 0.003549653112303053
$ws.Range("J3").Value = 0.003549653112303053
$ws.Range("O3").Value = 0.5596266124066729
$ws.Range("P3").Value = 0.5596266124066729
$ws.Range("Q3").Value = 0.03673660891722222
$ws.Range("R3").Value = 0.330629480255
$ws.Range("S3").Value = 0.001986480346456961
$ws.Range("T3").Value = 0.001986480346456961

$ws.Range("G4").Value = 0.003420333333333333
$ws.Range("H4").Value = 0.010261
$ws.Range("I4").Value = 0.003549653112303053
$ws.Range("J4").Value = 0.003549653112303053
$ws.Range("M4").Value = 7.027161333333335
$ws.Range("O4").Value = 0.3661403994706553
$ws.Range("P4").Value = 0.3661403994706552
$ws.Range("Q4").Value = 0.02403523414711111
$ws.Range("R4").Value = 0.216317107324
$ws.Range("S4").Value = 0.001299671408520895
$ws.Range("T4").Value = 0.001299671408520895

$ws.Range("I5").Value = 0.3907064193682856
$ws.Range("J5").Value = 0.3907064193682855
$ws.Range("M5").Value = 1.424719
$ws.Range("N5").Value = 4.274157
$ws.Range("O5").Value = 0.07423298812267187
$ws.Range("P5").Value = 0.07423298812267187
$ws.Range("Q5").Value = 0.5363672862743333
$ws.Range("R5").Value = 4.827305576469
$ws.Range("S5").Value = 0.0290033049884176
$ws.Range("T5").Value = 0.02900330498841759

$ws.Range("I6").Value = 0.3907064193682856
$ws.Range("J6").Value = 0.3907064193682855
$ws.Range("O6").Value = 0.5596266124066729
$ws.Range("P6").Value = 0.5596266124066729
$ws.Range("S6").Value = 0.2186497099166146
$ws.Range("T6").Value = 0.2186497099166145

$ws.Range("I7").Value = 0.3907064193682856
$ws.Range("J7").Value = 0.3907064193682855
$ws.Range("M7").Value = 7.027161333333335
$ws.Range("O7").Value = 0.3661403994706553
$ws.Range("P7").Value = 0.3661403994706552
$ws.Range("R7").Value = 23.80978641482801
$ws.Range("S7").Value = 0.1430534044632535
$ws.Range("T7").Value = 0.1430534044632534

$ws.Range("I8").Value = 0.6057439275194114
$ws.Range("J8").Value = 0.6057439275194113
$ws.Range("M8").Value = 1.424719
$ws.Range("N8").Value = 4.274157
$ws.Range("O8").Value = 0.07423298812267187
$ws.Range("P8").Value = 0.07423298812267187
$ws.Range("Q8").Value = 0.8315738121376667
$ws.Range("R8").Value = 7.484164309239
$ws.Range("S8").Value = 0.04496618177692908
$ws.Range("T8").Value = 0.04496618177692907

$ws.Range("I9").Value = 0.6057439275194114
$ws.Range("J9").Value = 0.6057439275194113
$ws.Range("O9").Value = 0.5596266124066729
$ws.Range("P9").Value = 0.5596266124066729
$ws.Range("S9").Value = 0.3389904221436014
$ws.Range("T9").Value = 0.3389904221436014

$ws.Range("I10").Value = 0.6057439275194114
$ws.Range("J10").Value = 0.6057439275194113
$ws.Range("M10").Value = 7.027161333333335
$ws.Range("O10").Value = 0.3661403994706553
$ws.Range("P10").Value = 0.3661403994706552
$ws.Range("R10").Value = 36.91424768406801
$ws.Range("S10").Value = 0.221787323598881
$ws.Range("T10").Value = 0.2217873235988809
